$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '37.839.23'
$ws.Range('D2').Style = $ws.Range('C2').Style
$ws.Range('E2').Value = '  -1.21%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.036.63'
$ws.Range('D3').Style = $ws.Range('C3').Style
$ws.Range('E3').Value = '  -1.39%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '227.27'
$ws.Range('D5').Style = $ws.Range('C5').Style
$ws.Range('E5').Value = '  -1.62%  '
$ws.Range('E6').Value = '  -0.63%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '60.16'
$ws.Range('D7').Style = $ws.Range('C7').Style
$ws.Range('E7').Value = '  +3.24%  '
$ws.Range('E8').Value = '  +0.17%  '
$ws.Range('E9').Value = '  -0.61%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0819'
$ws.Range('D10').Style = $ws.Range('C10').Style
$ws.Range('E10').Value = '  +1.40%  '
$ws.Range('E11').Value = '  +0.14%  '
$ws.Range('E12').Value = '  -0.20%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.336.24'
$ws.Range('D13').Style = $ws.Range('C13').Style
$ws.Range('E13').Value = '  -1.20%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '21.07'
$ws.Range('D14').Style = $ws.Range('C14').Style
$ws.Range('E14').Value = '  +1.50%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.761'
$ws.Range('D15').Style = $ws.Range('C15').Style
$ws.Range('E15').Value = '  +0.51%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.20'
$ws.Range('D16').Style = $ws.Range('C16').Style
$ws.Range('E16').Value = '  -2.04%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.041.35'
$ws.Range('D17').Style = $ws.Range('C17').Style
$ws.Range('E17').Value = '  -0.78%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '37.780.87'
$ws.Range('D18').Style = $ws.Range('C18').Style
$ws.Range('E18').Value = '  -1.00%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.07'
$ws.Range('D19').Style = $ws.Range('C19').Style
$ws.Range('E19').Value = '  -2.07%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '69.81'
$ws.Range('D20').Style = $ws.Range('C20').Style
$ws.Range('E20').Value = '  -0.15%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0825'
$ws.Range('D21').Style = $ws.Range('C21').Style
$ws.Range('E21').Value = '  -1.14%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '225.04'
$ws.Range('D22').Style = $ws.Range('C22').Style
$ws.Range('E22').Value = '  -0.06%  '
$ws.Range('E23').Value = '  +0.00%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.39'
$ws.Range('D24').Style = $ws.Range('C24').Style
$ws.Range('E24').Value = '  -2.41%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.21'
$ws.Range('D25').Style = $ws.Range('C25').Style
$ws.Range('E25').Value = '  -1.86%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.28'
$ws.Range('D26').Style = $ws.Range('C26').Style
$ws.Range('E26').Value = '  -0.81%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '165.17'
$ws.Range('D27').Style = $ws.Range('C27').Style
$ws.Range('E27').Value = '  -0.65%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.95'
$ws.Range('D29').Style = $ws.Range('C29').Style
$ws.Range('E29').Value = '  -0.85%  '
$ws.Range('E30').Value = '  -6.36%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.120'
$ws.Range('D31').Style = $ws.Range('C31').Style
$ws.Range('E31').Value = '  +1.14%  '
$ws.Range('E32').Value = '  -2.99%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.07'
$ws.Range('D33').Style = $ws.Range('C33').Style
$ws.Range('E33').Value = '  +4.39%  '
$ws.Range('E34').Value = '  -2.65%  '
$ws.Range('E35').Value = '  -2.65%  '
$ws.Range('E36').Value = '  +4.95%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.26'
$ws.Range('D37').Style = $ws.Range('C37').Style
$ws.Range('E37').Value = '  -5.52%  '
$ws.Range('E38').Value = '  -2.81%  '
$ws.Range('E39').Value = '  +0.18%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.545.63'
$ws.Range('D40').Style = $ws.Range('C40').Style
$ws.Range('E40').Value = '  +4.04%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0217'
$ws.Range('D41').Style = $ws.Range('C41').Style
$ws.Range('E41').Value = '  -0.75%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '97.07'
$ws.Range('D42').Style = $ws.Range('C42').Style
$ws.Range('E42').Value = '  -1.42%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '16.85'
$ws.Range('D43').Style = $ws.Range('C43').Style
$ws.Range('E43').Value = '  -1.14%  '
$ws.Range('E44').Value = '  -0.74%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0924'
$ws.Range('D45').Style = $ws.Range('C45').Style
$ws.Range('E45').Value = '  -2.35%  '
$ws.Range('E46').Value = '  -1.86%  '
$ws.Range('E47').Value = '  -4.75%  '
$ws.Range('E48').Value = '  -1.57%  '
$ws.Range('E49').Value = '  -0.30%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.226.43'
$ws.Range('D51').Style = $ws.Range('C51').Style
$ws.Range('E51').Value = '  -1.09%  '
